$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("startup")

# Update cell A2 from "CasesTab" to "ParticipantsTab"
$ws.Range("A2").Value = "ParticipantsTab"

# Update the selected/active cell to A2 (was B3)
$ws.Activate()
$ws.Range("A2").Select()
